$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 4.3
$ws.Range("G2").Value = 5.4
$ws.Range("H2").Value = 1.75
$ws.Range("I2").Value = 1.96
$ws.Range("J2").Value = 3.65
$ws.Range("K2").Value = 4.4
$ws.Range("L2").Value = 1.36
$ws.Range("N2").Value = 3.7
$ws.Range("O2").Value = 1.28
$ws.Range("P2").Value = 1.94
$ws.Range("Q2").Value = 1.81
$ws.Range("R2").Value = 1.36
$ws.Range("S2").Value = 2.86
$ws.Range("T2").Value = 1.78
$ws.Range("U2").Value = 2
$ws.Range("V2").Value = 2.04
$ws.Range("W2").Value = 1.23
$ws.Range("X2").Value = 17
$ws.Range("Y2").Value = 10
$ws.Range("Z2").Value = 12
$ws.Range("AA2").Value = 22
$ws.Range("AB2").Value = 18.5
$ws.Range("AC2").Value = 9.4
$ws.Range("AF2").Value = 40
$ws.Range("AG2").Value = 20
$ws.Range("AH2").Value = 21
$ws.Range("AI2").Value = 38
$ws.Range("AJ2").Value = 140
$ws.Range("AL2").Value = 260
$ws.Range("AM2").Value = 580
$ws.Range("AN2").Value = 310
$ws.Range("H3").Value = 1.1
$ws.Range("W3").Value = 1.03
$ws.Range("F5").Value = 1.74
$ws.Range("G5").Value = 1.96
$ws.Range("H5").Value = 4.3
$ws.Range("I5").Value = 6.4
$ws.Range("J5").Value = 3.35
$ws.Range("K5").Value = 4.1
$ws.Range("L5").Value = 1.43
$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 3.1
$ws.Range("O5").Value = 1.36
$ws.Range("P5").Value = 1.73
$ws.Range("Q5").Value = 2.06
$ws.Range("R5").Value = 1.28
$ws.Range("S5").Value = 3.45
$ws.Range("T5").Value = 1.9
$ws.Range("U5").Value = 1.84
$ws.Range("V5").Value = 1.2
$ws.Range("W5").Value = 2.04
$ws.Range("X5").Value = 90
$ws.Range("AB5").Value = 29
$ws.Range("AC5").Value = 14
$ws.Range("AF5").Value = 40
$ws.Range("AG5").Value = 40
$ws.Range("AJ5").Value = 900
$ws.Range("AK5").Value = 120
$ws.Range("AN5").Value = 65
$ws.Range("F6").Value = 2.14
$ws.Range("G6").Value = 2.56
$ws.Range("H6").Value = 3.05
$ws.Range("I6").Value = 3.9
$ws.Range("K6").Value = 4
$ws.Range("N6").Value = 3.3
$ws.Range("P6").Value = 1.79
$ws.Range("Q6").Value = 1.9
$ws.Range("S6").Value = 3.3
$ws.Range("T6").Value = 1.83
$ws.Range("V6").Value = 1.35
$ws.Range("W6").Value = 1.64
$ws.Range("X6").Value = 15
$ws.Range("AH6").Value = 19.5
$ws.Range("AJ6").Value = 34
$ws.Range("F7").Value = 5.8
$ws.Range("G7").Value = 9.800000000000001
$ws.Range("H7").Value = 1.49
$ws.Range("I7").Value = 1.65
$ws.Range("J7").Value = 3.6
$ws.Range("K7").Value = 5.6
$ws.Range("M7").Value = 1.06
$ws.Range("N7").Value = 3.3
$ws.Range("O7").Value = 1.29
$ws.Range("P7").Value = 1.9
$ws.Range("Q7").Value = 1.84
$ws.Range("R7").Value = 1.35
$ws.Range("S7").Value = 2.92
$ws.Range("T7").Value = 1.94
$ws.Range("U7").Value = 1.81
$ws.Range("V7").Value = 2.54
$ws.Range("W7").Value = 1.13
$ws.Range("Y7").Value = 29
$ws.Range("Z7").Value = 18.5
$ws.Range("AC7").Value = 42
$ws.Range("AD7").Value = 40
$ws.Range("AO7").Value = 29
$ws.Range("F8").Value = 1.73
$ws.Range("G8").Value = 1.88
$ws.Range("H8").Value = 5.6
$ws.Range("J8").Value = 3.25
$ws.Range("K8").Value = 4.2
$ws.Range("L8").Value = 1.45
$ws.Range("M8").Value = 1.08
$ws.Range("N8").Value = 2.9
$ws.Range("O8").Value = 1.39
$ws.Range("P8").Value = 1.66
$ws.Range("Q8").Value = 2.1
$ws.Range("R8").Value = 1.24
$ws.Range("S8").Value = 4.4
$ws.Range("T8").Value = 1.95
$ws.Range("U8").Value = 1.77
$ws.Range("V8").Value = 1.19
$ws.Range("W8").Value = 2.12
$ws.Range("X8").Value = 13
$ws.Range("Y8").Value = 17
$ws.Range("AB8").Value = 7.8
$ws.Range("AC8").Value = 9
$ws.Range("AD8").Value = 22
$ws.Range("AG8").Value = 11.5
$ws.Range("AH8").Value = 25
$ws.Range("AJ8").Value = 23
$ws.Range("AK8").Value = 25
$ws.Range("AN8").Value = 19
$ws.Range("J9").Value = 3.25
$ws.Range("V9").Value = 1.43
$ws.Range("AJ9").Value = 900
$ws.Range("F10").Value = 1.04
$ws.Range("V10").Value = 1.3
$ws.Range("H11").Value = 5.7
$ws.Range("I11").Value = 5.9
$ws.Range("Q11").Value = 1.48
$ws.Range("R11").Value = 1.76
$ws.Range("S11").Value = 2.18
$ws.Range("W11").Value = 2.6
$ws.Range("Y11").Value = 75
$ws.Range("AL11").Value = 46
$ws.Range("AM11").Value = 70
$ws.Range("AO11").Value = 42
$ws.Range("K12").Value = 3.75
$ws.Range("AK12").Value = 30
$ws.Range("AL12").Value = 150
$ws.Range("G13").Value = 3.6
$ws.Range("H13").Value = 2.1
$ws.Range("I13").Value = 2.2
$ws.Range("J13").Value = 3.9
$ws.Range("K13").Value = 4.1
$ws.Range("M13").Value = 1.05
$ws.Range("N13").Value = 4.9
$ws.Range("O13").Value = 1.22
$ws.Range("P13").Value = 2.46
$ws.Range("Q13").Value = 1.67
$ws.Range("R13").Value = 1.55
$ws.Range("S13").Value = 2.66
$ws.Range("T13").Value = 1.61
$ws.Range("U13").Value = 2.46
$ws.Range("V13").Value = 1.84
$ws.Range("W13").Value = 1.38
$ws.Range("X13").Value = 21
$ws.Range("Y13").Value = 13
$ws.Range("Z13").Value = 15.5
$ws.Range("AB13").Value = 18.5
$ws.Range("AC13").Value = 9
$ws.Range("AE13").Value = 20
$ws.Range("AF13").Value = 29
$ws.Range("AI13").Value = 29
$ws.Range("AK13").Value = 36
$ws.Range("AL13").Value = 40
$ws.Range("AM13").Value = 260
$ws.Range("AN13").Value = 26
$ws.Range("F14").Value = 6.2
$ws.Range("G14").Value = 6.6
$ws.Range("H14").Value = 1.61
$ws.Range("I14").Value = 1.62
$ws.Range("J14").Value = 4.4
$ws.Range("K14").Value = 4.5
$ws.Range("L14").Value = 1.33
$ws.Range("M14").Value = 1.05
$ws.Range("N14").Value = 4.4
$ws.Range("O14").Value = 1.26
$ws.Range("P14").Value = 2.12
$ws.Range("Q14").Value = 1.79
$ws.Range("R14").Value = 1.45
$ws.Range("S14").Value = 2.94
$ws.Range("T14").Value = 1.84
$ws.Range("U14").Value = 2.08
$ws.Range("V14").Value = 2.6
$ws.Range("W14").Value = 1.17
$ws.Range("X14").Value = 18
$ws.Range("Y14").Value = 9.6
$ws.Range("Z14").Value = 10
$ws.Range("AA14").Value = 15.5
$ws.Range("AB14").Value = 25
$ws.Range("AC14").Value = 9.4
$ws.Range("AD14").Value = 10
$ws.Range("AE14").Value = 16.5
$ws.Range("AG14").Value = 23
$ws.Range("AH14").Value = 21
$ws.Range("AI14").Value = 34
$ws.Range("AJ14").Value = 210
$ws.Range("AK14").Value = 85
$ws.Range("AL14").Value = 85
$ws.Range("AM14").Value = 330
$ws.Range("AN14").Value = 1000
$ws.Range("AO14").Value = 8.199999999999999
$ws.Range("P15").Value = 1.25
